$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "others"
$ws.Range("J2").Value = 2000

$ws.Range("O4").Select()
